$wb = $excel.ActiveWorkbook

# Work on the "Repayment Schedule" sheet: insert a new blank column before
# column N (pushing the old N/O/P columns to O/P/Q).
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Columns("N").Insert() | Out-Null

# Make "Repayment Schedule" the active sheet/tab and update its selection.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("S6").Select() | Out-Null
